$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Tuesday, Jan 10 - W92066 arrival
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Tuesday, Jan 10"
$ws.Range("C13").Value = "9:35 AM"
$ws.Range("D13").Value = "W92066"
$ws.Range("E13").Value = "London"
$ws.Range("F13").Value = "(LTN)"
$ws.Range("G13").Value = "Wizz Air "
$ws.Range("H13").Value = "A21N"
$ws.Range("I13").Value = "(G-WUKM)"
$ws.Range("J13").Value = "9:12 AM"
$ws.Range("K13").Font.Size = 11
$ws.Range("L13").Value = "0 hours, -23 minutes"
$ws.Range("M13").Font.Size = 11

# Row 14: Tuesday, Jan 10 - W62080 arrival
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Tuesday, Jan 10"
$ws.Range("C14").Value = "10:30 AM"
$ws.Range("D14").Value = "W62080"
$ws.Range("E14").Value = "Oslo"
$ws.Range("F14").Value = "(TRF)"
$ws.Range("G14").Value = "Wizz Air "
$ws.Range("H14").Value = "A321"
$ws.Range("I14").Value = "(HA-LTB)"
$ws.Range("J14").Value = "10:27 AM"
$ws.Range("K14").Font.Size = 11
$ws.Range("L14").Value = "0 hours, -3 minutes"
$ws.Range("M14").Font.Size = 11
